$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1): ID | User | Exchange ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "User"
$ws.Range("C1").Value = "Exchange"

# --- New data row (row 2): BFA.BB01.Dev | b.fargeout@outlook.com | ByBit ---
$ws.Range("A2").Value = "BFA.BB01.Dev"
$ws.Range("B2").Value = "b.fargeout@outlook.com"
$ws.Range("C2").Value = "ByBit"

# --- Row height for the header row ---
$ws.Rows.Item(1).RowHeight = 21

# --- Column widths (A, B, C) ---
$ws.Columns.Item(1).ColumnWidth = 11.833333333333334
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666

# --- Clear the old "vertical centered" alignment that used to live on the header style ---
$ws.Range("A1:C1").VerticalAlignment = -4107

# --- Both rows now use a Text ("@") number format ---
$ws.Range("A1:C2").NumberFormat = "@"
